# Daily attendance processing - 2025-12-09 10:31:48
#
# 1) Attendance sheet: 3 new late-scan records (Session 2, B1-3 group,
#    08/12/2025 12:30:00) are inserted at row 435, pushing every existing
#    record down by 3 rows (481 -> 484 total rows).
# 2) Summary sheet: the 3 affected students (211410, 211439, 211446) get
#    their rolled-up attendance counters refreshed to reflect the extra
#    session.
# 3) The Attendance sheet's AutoFilter range / _FilterDatabase defined name
#    are extended to cover the new last row (K484).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Attendance sheet - insert 3 new rows at the top of the B1-6 block
# ---------------------------------------------------------------------
$att = $wb.Worksheets.Item("Attendance")

$att.Range("A435:A437").EntireRow.Insert()

$newRows = @(
    @("211439", "عبد الرحمن سامح عبد العزيز منصور", "Year 5", "B1-3", "211439@med.asu.edu.eg", "GENERAL SURGERY", "2", "GENERAL SURGERY", "08/12/2025", "12:30:00", "B1-3"),
    @("211446", "ساره عبد الله محمد كمال عبد العزيز", "Year 5", "B1-3", "211446@med.asu.edu.eg", "GENERAL SURGERY", "2", "GENERAL SURGERY", "08/12/2025", "12:30:00", "B1-3"),
    @("211410", "جون مجدى ميخائيل سدراك", "Year 5", "B1-3", "211410@med.asu.edu.eg", "GENERAL SURGERY", "2", "GENERAL SURGERY", "08/12/2025", "12:30:00", "B1-3")
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 435 + $i
    $row = $newRows[$i]
    $att.Cells.Item($r, 1).Value = "'" + $row[0]
    $att.Cells.Item($r, 2).Value = $row[1]
    $att.Cells.Item($r, 3).Value = $row[2]
    $att.Cells.Item($r, 4).Value = $row[3]
    $att.Cells.Item($r, 5).Value = $row[4]
    $att.Cells.Item($r, 6).Value = $row[5]
    $att.Cells.Item($r, 7).Value = "'" + $row[6]
    $att.Cells.Item($r, 8).Value = $row[7]
    $att.Cells.Item($r, 9).Value = "'" + $row[8]
    $att.Cells.Item($r, 10).Value = "'" + $row[9]
    $att.Cells.Item($r, 11).Value = $row[10]
}

# Refresh the AutoFilter range and the _FilterDatabase defined name so
# they cover the new A1:K484 extent.
$att.AutoFilterMode = $false
$att.Range("A1:K484").AutoFilter()

$fdb = $wb.Names.Item("Attendance!_FilterDatabase")
$fdb.RefersTo = "='Attendance'!`$A`$1:`$K`$484"

# ---------------------------------------------------------------------
# 2. Summary sheet - refresh rollups for the 3 affected students
# ---------------------------------------------------------------------
$sum = $wb.Worksheets.Item("Summary")

$summaryRows = @(140, 159, 164)
foreach ($r in $summaryRows) {
    $sum.Cells.Item($r, 7).Value = "'13.3%"   # G: Percentage
    $sum.Cells.Item($r, 9).Value = 10          # I: Sessions Needed
    $sum.Cells.Item($r, 14).Value = 2          # N: Total Attended
    $sum.Cells.Item($r, 15).Value = 0          # O: Total Missed
    $sum.Cells.Item($r, 17).Value = 2          # Q: Attended GENERAL SURGERY (Total)
    $sum.Cells.Item($r, 19).Value = 1          # S: GENERAL SURGERY Session 2
}
